# Apply updated evaluation metric values across the three sheets:
# Summary, Classification Report, Confusion Matrix

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 0.5
$wsSummary.Range("C2").Value = 0.5
$wsSummary.Range("D2").Value = 1
$wsSummary.Range("E2").Value = 0.6666666666666666
$wsSummary.Range("F2").Value = 0.8333333333333334
$wsSummary.Range("G2").Value = 0.9629629629629629
$wsSummary.Range("H2").Value = 0.7963605885901051
$wsSummary.Range("I2").Value = 534
$wsSummary.Range("J2").Value = 534
$wsSummary.Range("K2").Value = 0
$wsSummary.Range("L2").Value = 0

# ---------------------------------------------------------------
# Sheet: Classification Report
# ---------------------------------------------------------------
$wsClass = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$wsClass.Range("B2").Value = 0
$wsClass.Range("C2").Value = 0
$wsClass.Range("D2").Value = 0

# Row 3 - class "1"
$wsClass.Range("B3").Value = 0.5
$wsClass.Range("C3").Value = 1
$wsClass.Range("D3").Value = 0.6666666666666666

# Row 4 - accuracy
$wsClass.Range("B4").Value = 0.5
$wsClass.Range("C4").Value = 0.5
$wsClass.Range("D4").Value = 0.5
$wsClass.Range("E4").Value = 0.5

# Row 5 - macro avg
$wsClass.Range("B5").Value = 0.25
$wsClass.Range("C5").Value = 0.5
$wsClass.Range("D5").Value = 0.3333333333333333

# Row 6 - weighted avg
$wsClass.Range("B6").Value = 0.25
$wsClass.Range("C6").Value = 0.5
$wsClass.Range("D6").Value = 0.3333333333333333

# ---------------------------------------------------------------
# Sheet: Confusion Matrix
# ---------------------------------------------------------------
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$wsConf.Range("B2").Value = 0
$wsConf.Range("C2").Value = 534

# Row 3 - Actual 1
$wsConf.Range("B3").Value = 0
$wsConf.Range("C3").Value = 534
